$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1353624.5
$ws.Range("I15").Value = 1353624.5
$ws.Range("K15").Value = 4060873.5
$ws.Range("M15").Value = -4060704.5
$ws.Range("H40").Value = 2139.1304
$ws.Range("I40").Value = 1512.5883
$ws.Range("J40").Value = 3914.3333
$ws.Range("K40").Value = 1512.5883
$ws.Range("L40").Value = 3914.3333
$ws.Range("M40").Value = -1337.5883
$ws.Range("N40").Value = -4264.3333
$ws.Range("H64").Value = 3287.7058
$ws.Range("I64").Value = 2820.3
$ws.Range("J64").Value = 3482.4583
$ws.Range("K64").Value = 2820.3
$ws.Range("L64").Value = 3482.4583
$ws.Range("M64").Value = -2572.3
$ws.Range("N64").Value = -3978.4583
$ws.Range("H67").Value = 3287.7058
$ws.Range("I67").Value = 2820.3
$ws.Range("J67").Value = 3482.4583
$ws.Range("K67").Value = 2820.3
$ws.Range("L67").Value = 3482.4583
$ws.Range("M67").Value = -1962.3
$ws.Range("N67").Value = -5198.4583
$ws.Range("H112").Value = 2827.2666
$ws.Range("I112").Value = 495
$ws.Range("J112").Value = 3186.077
$ws.Range("K112").Value = 1485
$ws.Range("L112").Value = 9558.231
$ws.Range("M112").Value = -377
$ws.Range("N112").Value = -11774.231
$ws.Range("H113").Value = 3677.5
$ws.Range("I113").Value = 2695.3333
$ws.Range("J113").Value = 4856.1
$ws.Range("K113").Value = 2695.3333
$ws.Range("L113").Value = 4856.1
$ws.Range("M113").Value = 558.6667000000002
$ws.Range("N113").Value = -11364.1
$ws.Range("H116").Value = 6963.684
$ws.Range("I116").Value = 4037.1428
$ws.Range("J116").Value = 15158
$ws.Range("K116").Value = 4037.1428
$ws.Range("L116").Value = 15158
$ws.Range("M116").Value = -595.1428000000001
$ws.Range("N116").Value = -22042
$ws.Range("H132").Value = 6851.1465
$ws.Range("I132").Value = 5153.0312
$ws.Range("J132").Value = 12888.889
$ws.Range("K132").Value = 15459.0936
$ws.Range("L132").Value = 38666.667
$ws.Range("M132").Value = -12929.0936
$ws.Range("N132").Value = -43726.667
$ws.Range("H137").Value = 1417.8462
$ws.Range("I137").Value = 1875.3334
$ws.Range("J137").Value = 1025.7142
$ws.Range("K137").Value = 5626.0002
$ws.Range("L137").Value = 3077.1426
$ws.Range("M137").Value = -3076.0002
$ws.Range("N137").Value = -8177.142599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1026.1364
$ws.Range("I2").Value = 450.92307
$ws.Range("J2").Value = 1857
$ws.Range("K2").Value = 450.92307
$ws.Range("L2").Value = 1857
$ws.Range("M2").Value = -337.92307
$ws.Range("N2").Value = -2083
$ws.Range("H32").Value = 7431.5625
$ws.Range("I32").Value = 6567
$ws.Range("K32").Value = 6567
$ws.Range("M32").Value = -6280
$ws.Range("H61").Value = 7203.1665
$ws.Range("I61").Value = 8162.4
$ws.Range("J61").Value = 2407
$ws.Range("K61").Value = 8162.4
$ws.Range("L61").Value = 2407
$ws.Range("M61").Value = -7950.4
$ws.Range("N61").Value = -2831
$ws.Range("H63").Value = 6562.857
$ws.Range("I63").Value = 3500
$ws.Range("J63").Value = 10646.667
$ws.Range("K63").Value = 3500
$ws.Range("L63").Value = 10646.667
$ws.Range("M63").Value = -2814
$ws.Range("N63").Value = -12018.667
$ws.Range("H66").Value = 6562.857
$ws.Range("I66").Value = 3500
$ws.Range("J66").Value = 10646.667
$ws.Range("K66").Value = 17500
$ws.Range("L66").Value = 53233.335
$ws.Range("M66").Value = -14068
$ws.Range("N66").Value = -60097.335
$ws.Range("H74").Value = 1717.2354
$ws.Range("I74").Value = 1699.3636
$ws.Range("J74").Value = 1750
$ws.Range("K74").Value = 1699.3636
$ws.Range("L74").Value = 1750
$ws.Range("M74").Value = -825.3635999999999
$ws.Range("N74").Value = -3498
$ws.Range("H77").Value = 1717.2354
$ws.Range("I77").Value = 1699.3636
$ws.Range("J77").Value = 1750
$ws.Range("K77").Value = 8496.817999999999
$ws.Range("L77").Value = 8750
$ws.Range("M77").Value = -4128.817999999999
$ws.Range("N77").Value = -17486
$ws.Range("H116").Value = 1026.1364
$ws.Range("I116").Value = 450.92307
$ws.Range("J116").Value = 1857
$ws.Range("K116").Value = 450.92307
$ws.Range("L116").Value = 1857
$ws.Range("M116").Value = 1843.07693
$ws.Range("N116").Value = -6445
$ws.Range("H132").Value = 638287.4399999999
$ws.Range("I132").Value = 1222913.6
$ws.Range("J132").Value = 4942.4443
$ws.Range("K132").Value = 3668740.8
$ws.Range("L132").Value = 14827.3329
$ws.Range("M132").Value = -3666210.8
$ws.Range("N132").Value = -19887.3329
$ws.Range("H136").Value = 7203.1665
$ws.Range("I136").Value = 8162.4
$ws.Range("J136").Value = 2407
$ws.Range("K136").Value = 24487.2
$ws.Range("L136").Value = 7221
$ws.Range("M136").Value = -21937.2
$ws.Range("N136").Value = -12321

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1026.1364
$ws.Range("I3").Value = 450.92307
$ws.Range("J3").Value = 1857
$ws.Range("K3").Value = 450.92307
$ws.Range("L3").Value = 1857
$ws.Range("M3").Value = -336.92307
$ws.Range("N3").Value = -2085
$ws.Range("H134").Value = 5154.3403
$ws.Range("I134").Value = 2136.1738
$ws.Range("J134").Value = 8046.75
$ws.Range("K134").Value = 6408.5214
$ws.Range("L134").Value = 24140.25
$ws.Range("M134").Value = -3873.5214
$ws.Range("N134").Value = -29210.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7790.1177
$ws.Range("I58").Value = 3935.3333
$ws.Range("J58").Value = 12126.75
$ws.Range("K58").Value = 3935.3333
$ws.Range("L58").Value = 12126.75
$ws.Range("M58").Value = -3732.3333
$ws.Range("N58").Value = -12532.75
$ws.Range("H132").Value = 2711.9429
$ws.Range("I132").Value = 2070.4
$ws.Range("K132").Value = 6211.200000000001
$ws.Range("M132").Value = -3681.200000000001
$ws.Range("H134").Value = 2516.182
$ws.Range("I134").Value = 1684.1428
$ws.Range("J134").Value = 3972.25
$ws.Range("K134").Value = 5052.428400000001
$ws.Range("L134").Value = 11916.75
$ws.Range("M134").Value = -2517.428400000001
$ws.Range("N134").Value = -16986.75
$ws.Range("H136").Value = 7790.1177
$ws.Range("I136").Value = 3935.3333
$ws.Range("J136").Value = 12126.75
$ws.Range("K136").Value = 11805.9999
$ws.Range("L136").Value = 36380.25
$ws.Range("M136").Value = -9255.999899999999
$ws.Range("N136").Value = -41480.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 971048.5600000001
$ws.Range("I132").Value = 1603998.5
$ws.Range("J132").Value = 3007.4119
$ws.Range("K132").Value = 4811995.5
$ws.Range("L132").Value = 9022.235700000001
$ws.Range("M132").Value = -4809465.5
$ws.Range("N132").Value = -14082.2357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1667446.6
$ws.Range("I46").Value = 630
$ws.Range("J46").Value = 3334263.2
$ws.Range("K46").Value = 630
$ws.Range("L46").Value = 3334263.2
$ws.Range("M46").Value = -442
$ws.Range("N46").Value = -3334639.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2210.4707
$ws.Range("I122").Value = 1405.9166
$ws.Range("K122").Value = 4217.7498
$ws.Range("M122").Value = -1767.7498
